$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 144445120
$ws.Range("J43").Value = 2000
$ws.Range("L43").Value = 2000
$ws.Range("N43").Value = -2138
$ws.Range("H51").Value = 100212000
$ws.Range("J51").Value = 250010000
$ws.Range("L51").Value = 250010000
$ws.Range("N51").Value = -250010968
$ws.Range("H100").Value = 2048.6
$ws.Range("I100").Value = 1837.7858
$ws.Range("K100").Value = 1837.7858
$ws.Range("M100").Value = -1296.7858
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H112").Value = 4565.1055
$ws.Range("J112").Value = 4996.0625
$ws.Range("L112").Value = 14988.1875
$ws.Range("N112").Value = -17204.1875
$ws.Range("H116").Value = 8078.8
$ws.Range("I116").Value = 7166.3335
$ws.Range("K116").Value = 7166.3335
$ws.Range("M116").Value = -3724.3335
$ws.Range("H138").Value = 4977.2964
$ws.Range("I138").Value = 6460.143
$ws.Range("K138").Value = 19380.429
$ws.Range("M138").Value = -14240.429
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7801.125
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 7801.125
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 7801.125
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -8375.125
$ws.Range("H122").Value = 1908.8334
$ws.Range("I122").Value = 1946
$ws.Range("K122").Value = 5838
$ws.Range("M122").Value = -3388
$ws.Range("H132").Value = 5552.1665
$ws.Range("I132").Value = 5463.6
$ws.Range("K132").Value = 16390.8
$ws.Range("M132").Value = -13860.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 6206.5
$ws.Range("I94").Value = 5823.875
$ws.Range("J94").Value = 8502.25
$ws.Range("K94").Value = 5823.875
$ws.Range("L94").Value = 8502.25
$ws.Range("M94").Value = -5372.875
$ws.Range("N94").Value = -9404.25
$ws.Range("H99").Value = 10478.214
$ws.Range("I99").Value = 18813.572
$ws.Range("J99").Value = 2142.8572
$ws.Range("K99").Value = 18813.572
$ws.Range("L99").Value = 2142.8572
$ws.Range("M99").Value = -17315.572
$ws.Range("N99").Value = -5138.8572
$ws.Range("H107").Value = 7045.9287
$ws.Range("I107").Value = 7850.8857
$ws.Range("K107").Value = 7850.8857
$ws.Range("M107").Value = -5930.8857
$ws.Range("H134").Value = 28127238
$ws.Range("J134").Value = 60002068
$ws.Range("L134").Value = 180006204
$ws.Range("N134").Value = -180011274
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 863.5333000000001
$ws.Range("J22").Value = 564.7143
$ws.Range("L22").Value = 564.7143
$ws.Range("N22").Value = -1264.7143
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H62").Value = 1518.4
$ws.Range("J62").Value = 1196
$ws.Range("L62").Value = 1196
$ws.Range("N62").Value = -2444
$ws.Range("H65").Value = 1518.4
$ws.Range("J65").Value = 1196
$ws.Range("L65").Value = 5980
$ws.Range("N65").Value = -12220
$ws.Range("H86").Value = 14127
$ws.Range("I86").Value = 8732.200000000001
$ws.Range("K86").Value = 8732.200000000001
$ws.Range("M86").Value = -7609.200000000001
$ws.Range("H89").Value = 14127
$ws.Range("I89").Value = 8732.200000000001
$ws.Range("K89").Value = 43661
$ws.Range("M89").Value = -38045
$ws.Range("H105").Value = 2626.0908
$ws.Range("I105").Value = 1954.625
$ws.Range("K105").Value = 1954.625
$ws.Range("M105").Value = -207.625
$ws.Range("H107").Value = 1781.375
$ws.Range("J107").Value = 10000
$ws.Range("L107").Value = 10000
$ws.Range("N107").Value = -13840
$ws.Range("H122").Value = 2570.4666
$ws.Range("I122").Value = 2700.6924
$ws.Range("J122").Value = 1724
$ws.Range("K122").Value = 8102.0772
$ws.Range("L122").Value = 5172
$ws.Range("M122").Value = -5652.0772
$ws.Range("N122").Value = -10072
$ws.Range("H132").Value = 24521.979
$ws.Range("I132").Value = 32388.334
$ws.Range("K132").Value = 97165.00199999999
$ws.Range("M132").Value = -94635.00199999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2070.8572
$ws.Range("I68").Value = 2074.75
$ws.Range("K68").Value = 6224.25
$ws.Range("M68").Value = -5413.25
$ws.Range("H71").Value = 2070.8572
$ws.Range("I71").Value = 2074.75
$ws.Range("K71").Value = 18672.75
$ws.Range("M71").Value = -14616.75
$ws.Range("H107").Value = 111112180
$ws.Range("J107").Value = 111112180
$ws.Range("L107").Value = 333336540
$ws.Range("N107").Value = -333340380
$ws.Range("H140").Value = 13890707
$ws.Range("I140").Value = 14494216
$ws.Range("J140").Value = 9999
$ws.Range("K140").Value = 43482648
$ws.Range("L140").Value = 29997
$ws.Range("M140").Value = -43477468
$ws.Range("N140").Value = -40357
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 83535000
$ws.Range("I80").Value = 368332
$ws.Range("J80").Value = 166701660
$ws.Range("K80").Value = 368332
$ws.Range("L80").Value = 166701660
$ws.Range("M80").Value = -367334
$ws.Range("N80").Value = -166703656
$ws.Range("H83").Value = 83535000
$ws.Range("I83").Value = 368332
$ws.Range("J83").Value = 166701660
$ws.Range("K83").Value = 1841660
$ws.Range("L83").Value = 833508300
$ws.Range("M83").Value = -1836668
$ws.Range("N83").Value = -833518284
$ws.Range("H107").Value = 100434.7
$ws.Range("J107").Value = 499.85715
$ws.Range("L107").Value = 499.85715
$ws.Range("N107").Value = -4339.85715
$ws.Range("H132").Value = 528333.6
$ws.Range("I132").Value = 6725
$ws.Range("K132").Value = 20175
$ws.Range("M132").Value = -17645
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9907.333000000001
$ws.Range("J2").Value = 13111
$ws.Range("L2").Value = 13111
$ws.Range("N2").Value = -13335
$ws.Range("H95").Value = 50343.5
$ws.Range("J95").Value = 50343.5
$ws.Range("L95").Value = 50343.5
$ws.Range("N95").Value = -55835.5
$ws.Range("H100").Value = 2683.4167
$ws.Range("I100").Value = 2785.1
$ws.Range("J100").Value = 2175
$ws.Range("K100").Value = 2785.1
$ws.Range("L100").Value = 2175
$ws.Range("M100").Value = -2244.1
$ws.Range("N100").Value = -3257
$ws.Range("H122").Value = 3245.0435
$ws.Range("I122").Value = 2878.5312
$ws.Range("K122").Value = 8635.5936
$ws.Range("M122").Value = -6185.5936
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3823
$ws.Range("I62").Value = 3766.8333
$ws.Range("K62").Value = 3766.8333
$ws.Range("M62").Value = -3142.8333
$ws.Range("H65").Value = 3823
$ws.Range("I65").Value = 3766.8333
$ws.Range("K65").Value = 18834.1665
$ws.Range("M65").Value = -15714.1665
$ws.Range("H81").Value = 78756.766
$ws.Range("I81").Value = 2184.4
$ws.Range("K81").Value = 4368.8
$ws.Range("M81").Value = -3307.8
$ws.Range("H84").Value = 78756.766
$ws.Range("I84").Value = 2184.4
$ws.Range("K84").Value = 21844
$ws.Range("M84").Value = -16540
$ws.Range("H113").Value = 1084.3334
$ws.Range("I113").Value = 1066.75
$ws.Range("K113").Value = 3200.25
$ws.Range("M113").Value = -1030.25
$ws.Range("H132").Value = 2588.318
$ws.Range("I132").Value = 2293.4783
$ws.Range("J132").Value = 2911.238
$ws.Range("K132").Value = 6880.4349
$ws.Range("L132").Value = 8733.714
$ws.Range("M132").Value = -4350.4349
$ws.Range("N132").Value = -13793.714
